$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Mark the six team members' tasks for this week as completed.
$ws.Range("C43").Value = "已完成"
$ws.Range("C44").Value = "已完成"
$ws.Range("C45").Value = "已完成"
$ws.Range("C46").Value = "已完成"
$ws.Range("C47").Value = "已完成"
$ws.Range("C48").Value = "已完成"

# Fill in the stage summary (previously just a bare "总结：" placeholder).
$ws.Range("A49").Value = "总结：此阶段除了完成以上计划内容之外，小组成员还完成了群签到管理界面设计、群签到界面设计。同时，对之前的UI设计不合理的位置进行了更改，对用例规约进行修正，并将UI图置入用例规约说明书中。基本完成了UI设计，但有些细节处还未处理完全。"

# Leave the selection where the author last left it.
$ws.Range("E48").Select()
